$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K2").Value = 2.12
$ws.Range("Q2").Value = 1.78
$ws.Range("R2").Value = 1.98
$ws.Range("S2").Value = 1.37
$ws.Range("T2").Value = 2.85
$ws.Range("X2").Value = 12
$ws.Range("AE2").Value = 11.75
$ws.Range("AH2").Value = 10.75
$ws.Range("AI2").Value = 18.5
$ws.Range("AM2").Value = 30
$ws.Range("AQ2").Value = 40
$ws.Range("AS2").Value = 175
$ws.Range("AT2").Value = 2.85
